$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data. Column D (Price) values are kept as plain
# text (number format "@") since the source data stores prices as literal
# strings (e.g. "30.127.54", "7.380") rather than numeric values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.127.54'
$ws.Range("E2").Value = '  +5.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.915.79'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("E4").Value = '  -0.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.31'
$ws.Range("E5").Value = '  +4.40%  '
$ws.Range("E6").Value = '  -0.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5232'
$ws.Range("E7").Value = '  +2.95%  '
$ws.Range("E8").Value = '  +4.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08507'
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.94'
$ws.Range("E10").Value = '  +1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.122'
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.35'
$ws.Range("E12").Value = '  +9.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.465'
$ws.Range("E13").Value = '  +4.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.916.91'
$ws.Range("E14").Value = '  +2.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.380'
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.98'
$ws.Range("E17").Value = '  +4.04%  '
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06687'
$ws.Range("E19").Value = '  -0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.37'
$ws.Range("E20").Value = '  +3.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.009'
$ws.Range("E22").Value = '  +1.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.127.84'
$ws.Range("E23").Value = '  +5.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.31'
$ws.Range("E24").Value = '  +1.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.216'
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.133.04'
$ws.Range("E26").Value = '  +2.56%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.60'
$ws.Range("E27").Value = '  +2.38%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.13'
$ws.Range("E28").Value = '  +2.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.416'
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.83'
$ws.Range("E30").Value = '  +2.23%  '
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("E32").Value = '  +2.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.017'
$ws.Range("E33").Value = '  +4.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.642'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02490'
$ws.Range("E35").Value = '  +1.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06586'
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2207'
$ws.Range("E37").Value = '  +2.08%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.176'
$ws.Range("E38").Value = '  +2.71%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.225'
$ws.Range("E39").Value = '  +3.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.848'
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6528'
$ws.Range("E41").Value = '  +2.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.61'
$ws.Range("E42").Value = '  +4.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.241'
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6139'
$ws.Range("E44").Value = '  +2.37%  '
$ws.Range("E45").Value = '  +1.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.741'
$ws.Range("E46").Value = '  +1.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.080'
$ws.Range("E47").Value = '  +3.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.243'
$ws.Range("E48").Value = '  +2.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.27'
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("E50").Value = '  +4.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.60'
$ws.Range("E51").Value = '  +4.64%  '
